# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the rows
# whose underlying item prices changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 252224.75
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 336133
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 1008399
$ws.Range("M46").Value = -1381
$ws.Range("N46").Value = -1008637
$ws.Range("H58").Value = 1464.0741
$ws.Range("J58").Value = 3817.5
$ws.Range("L58").Value = 11452.5
$ws.Range("N58").Value = -11752.5
$ws.Range("H59").Value = 1700
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 252224.75
$ws.Range("I60").Value = 500
$ws.Range("J60").Value = 336133
$ws.Range("K60").Value = 1500
$ws.Range("L60").Value = 1008399
$ws.Range("M60").Value = -1016
$ws.Range("N60").Value = -1009367
$ws.Range("H137").Value = 7740.2
$ws.Range("I137").Value = 11166.131
$ws.Range("K137").Value = 33498.393
$ws.Range("M137").Value = -30948.393
$ws.Range("H141").Value = 4149.5854
$ws.Range("I141").Value = 3261.697
$ws.Range("K141").Value = 9785.091
$ws.Range("M141").Value = -4605.091
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2240
$ws.Range("I10").Value = 130
$ws.Range("J10").Value = 4350
$ws.Range("K10").Value = 130
$ws.Range("L10").Value = 4350
$ws.Range("M10").Value = 10
$ws.Range("N10").Value = -4630
$ws.Range("H20").Value = 3103.0908
$ws.Range("I20").Value = 1460.4445
$ws.Range("K20").Value = 1460.4445
$ws.Range("M20").Value = -1213.4445
$ws.Range("H64").Value = 5992.222
$ws.Range("J64").Value = 2746.4546
$ws.Range("L64").Value = 2746.4546
$ws.Range("N64").Value = -3196.4546
$ws.Range("H67").Value = 5992.222
$ws.Range("J67").Value = 2746.4546
$ws.Range("L67").Value = 2746.4546
$ws.Range("N67").Value = -4306.4546
$ws.Range("H86").Value = 8186.1577
$ws.Range("I86").Value = 5683
$ws.Range("K86").Value = 5683
$ws.Range("M86").Value = -4560
$ws.Range("H89").Value = 8186.1577
$ws.Range("I89").Value = 5683
$ws.Range("K89").Value = 28415
$ws.Range("M89").Value = -22799
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2240.6667
$ws.Range("J58").Value = 3888.4285
$ws.Range("L58").Value = 3888.4285
$ws.Range("N58").Value = -4294.4285
$ws.Range("H74").Value = 64219.668
$ws.Range("J74").Value = 81140.28999999999
$ws.Range("L74").Value = 81140.28999999999
$ws.Range("N74").Value = -82888.28999999999
$ws.Range("H77").Value = 64219.668
$ws.Range("J77").Value = 81140.28999999999
$ws.Range("L77").Value = 243420.87
$ws.Range("N77").Value = -252156.87
$ws.Range("H94").Value = 2320.8333
$ws.Range("I94").Value = 5628.25
$ws.Range("K94").Value = 5628.25
$ws.Range("M94").Value = -5177.25
$ws.Range("H132").Value = 16317.179
$ws.Range("I132").Value = 1232.0555
$ws.Range("K132").Value = 3696.1665
$ws.Range("M132").Value = -1166.1665
$ws.Range("H134").Value = 5965.579
$ws.Range("I134").Value = 3276.5833
$ws.Range("K134").Value = 9829.749899999999
$ws.Range("M134").Value = -7294.749899999999
$ws.Range("H136").Value = 2240.6667
$ws.Range("J136").Value = 3888.4285
$ws.Range("L136").Value = 11665.2855
$ws.Range("N136").Value = -16765.2855
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 366.66666
$ws.Range("J75").Value = 400
$ws.Range("L75").Value = 1200
$ws.Range("N75").Value = -3196
$ws.Range("H78").Value = 366.66666
$ws.Range("J78").Value = 400
$ws.Range("L78").Value = 3600
$ws.Range("N78").Value = -13584
$ws.Range("H80").Value = 43607.5
$ws.Range("J80").Value = 61580.855
$ws.Range("L80").Value = 184742.565
$ws.Range("N80").Value = -186614.565
$ws.Range("H83").Value = 43607.5
$ws.Range("J83").Value = 61580.855
$ws.Range("L83").Value = 554227.6950000001
$ws.Range("N83").Value = -563587.6950000001
$ws.Range("H98").Value = 1481.75
$ws.Range("J98").Value = 1350
$ws.Range("L98").Value = 4050
$ws.Range("N98").Value = -7046
$ws.Range("H107").Value = 606.2727
$ws.Range("I107").Value = 207.91667
$ws.Range("J107").Value = 755.65625
$ws.Range("K107").Value = 623.75001
$ws.Range("L107").Value = 2266.96875
$ws.Range("M107").Value = 1296.24999
$ws.Range("N107").Value = -6106.96875
$ws.Range("H117").Value = 1226.6364
$ws.Range("J117").Value = 352.33334
$ws.Range("L117").Value = 1057.00002
$ws.Range("N117").Value = -7941.000019999999
$ws.Range("H132").Value = 43575
$ws.Range("J132").Value = 73971.14
$ws.Range("L132").Value = 665740.26
$ws.Range("N132").Value = -670800.26
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9624
$ws.Range("I80").Value = 19000
$ws.Range("J80").Value = 3998.4
$ws.Range("K80").Value = 19000
$ws.Range("L80").Value = 3998.4
$ws.Range("M80").Value = -18002
$ws.Range("N80").Value = -5994.4
$ws.Range("H83").Value = 9624
$ws.Range("I83").Value = 19000
$ws.Range("J83").Value = 3998.4
$ws.Range("K83").Value = 95000
$ws.Range("L83").Value = 19992
$ws.Range("M83").Value = -90008
$ws.Range("N83").Value = -29976
$ws.Range("H102").Value = 10363.706
$ws.Range("I102").Value = 14107.546
$ws.Range("K102").Value = 14107.546
$ws.Range("M102").Value = -12485.546
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 14459.7
$ws.Range("I22").Value = 16387.125
$ws.Range("K22").Value = 16387.125
$ws.Range("M22").Value = -16092.125
$ws.Range("H27").Value = 14459.7
$ws.Range("I27").Value = 16387.125
$ws.Range("K27").Value = 16387.125
$ws.Range("M27").Value = -16280.125
$ws.Range("H40").Value = 33617.812
$ws.Range("I40").Value = 44899.5
$ws.Range("J40").Value = 14815
$ws.Range("K40").Value = 44899.5
$ws.Range("L40").Value = 14815
$ws.Range("M40").Value = -44763.5
$ws.Range("N40").Value = -15087
$ws.Range("H68").Value = 6869.8
$ws.Range("I68").Value = 2283
$ws.Range("K68").Value = 2283
$ws.Range("M68").Value = -1534
$ws.Range("H71").Value = 6869.8
$ws.Range("I71").Value = 2283
$ws.Range("K71").Value = 11415
$ws.Range("M71").Value = -7671
$ws.Range("H122").Value = 5691.4
$ws.Range("I122").Value = 5691.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17074.2
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14624.2
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 7777.1113
$ws.Range("I136").Value = 598
$ws.Range("J136").Value = 8199.412
$ws.Range("K136").Value = 1794
$ws.Range("L136").Value = 24598.236
$ws.Range("M136").Value = 756
$ws.Range("N136").Value = -29698.236
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 202649.69
$ws.Range("I62").Value = 635224.8
$ws.Range("K62").Value = 635224.8
$ws.Range("M62").Value = -634600.8
$ws.Range("H65").Value = 202649.69
$ws.Range("I65").Value = 635224.8
$ws.Range("K65").Value = 3176124
$ws.Range("M65").Value = -3173004
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H81").Value = 12491.5
$ws.Range("I81").Value = 13656.75
$ws.Range("J81").Value = 5500
$ws.Range("K81").Value = 27313.5
$ws.Range("L81").Value = 11000
$ws.Range("M81").Value = -26252.5
$ws.Range("N81").Value = -13122
$ws.Range("H84").Value = 12491.5
$ws.Range("I84").Value = 13656.75
$ws.Range("J84").Value = 5500
$ws.Range("K84").Value = 136567.5
$ws.Range("L84").Value = 55000
$ws.Range("M84").Value = -131263.5
$ws.Range("N84").Value = -65608
$ws.Range("H122").Value = 3224.913
$ws.Range("I122").Value = 3224.913
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9674.739
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7224.739
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 27852
$ws.Range("I126").Value = 41849.7
$ws.Range("K126").Value = 125549.1
$ws.Range("M126").Value = -123079.1
$ws.Range("H131").Value = 39650
$ws.Range("I131").Value = 39650
$ws.Range("K131").Value = 39650
$ws.Range("M131").Value = -34610
$ws.Range("H132").Value = 43932.223
$ws.Range("I132").Value = 115800
$ws.Range("J132").Value = 7998.3335
$ws.Range("K132").Value = 347400
$ws.Range("L132").Value = 23995.0005
$ws.Range("M132").Value = -344870
$ws.Range("N132").Value = -29055.0005
$ws.Range("H136").Value = 2000.0605
$ws.Range("I136").Value = 1240
$ws.Range("J136").Value = 4375.25
$ws.Range("K136").Value = 3720
$ws.Range("L136").Value = 13125.75
$ws.Range("M136").Value = -1170
$ws.Range("N136").Value = -18225.75
